# Actualización automática 2025-06-06 15:25:08
# Insert a new client row ("GOYA CASTRO CARLOS ALBERTO") in alphabetical
# order (row 26) on both the "VENTAS POR GRUPO" and "VENTA MENSUAL"
# sheets, shifting the existing rows down by one, and refresh the
# trailing summary row's "X de 53" -> "X de 54" counters.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO  (columns A:N, data rows 2-54, summary 55->56)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows("26:26").Insert()

$ws1.Cells.Item(26, 1).Value = "CASTRO ALCIVAR EDA MARIA"
$ws1.Cells.Item(26, 2).Value = "GOYA CASTRO CARLOS ALBERTO"
for ($col = 3; $col -le 14; $col++) {
    $ws1.Cells.Item(26, $col).Value = 0
}

for ($col = 3; $col -le 14; $col++) {
    $cell = $ws1.Cells.Item(56, $col)
    $text = $cell.Value()
    $cell.Value = $text.Replace("de 53", "de 54")
}

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL  (columns A:G, data rows 2-54, summary 55->56)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows("26:26").Insert()

$ws2.Cells.Item(26, 1).Value = "CASTRO ALCIVAR EDA MARIA"
$ws2.Cells.Item(26, 2).Value = "GOYA CASTRO CARLOS ALBERTO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(26, $col).Value = 0
}
